$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.226.79"
$ws.Range("E2").Value = "  -3.20%  "
$ws.Range("D3").Value = "'1.729.69"
$ws.Range("E3").Value = "  -3.73%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'322.39"
$ws.Range("E5").Value = "  -4.47%  "
$ws.Range("D7").Value = "'0.4215"
$ws.Range("E7").Value = "  -8.46%  "
$ws.Range("D8").Value = "'0.3571"
$ws.Range("E8").Value = "  -3.32%  "
$ws.Range("D9").Value = "'44.85"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").Value = "'0.07383"
$ws.Range("E10").Value = "  -3.73%  "
$ws.Range("D11").Value = "'1.105"
$ws.Range("E11").Value = "  -3.41%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").Value = "'21.44"
$ws.Range("E13").Value = "  -4.48%  "
$ws.Range("D14").Value = "'6.046"
$ws.Range("E14").Value = "  -4.62%  "
$ws.Range("D15").Value = "'7.108"
$ws.Range("E15").Value = "  -3.72%  "
$ws.Range("D16").Value = "'1.730.50"
$ws.Range("E16").Value = "  -3.55%  "
$ws.Range("D17").Value = "'0.00001055"
$ws.Range("E17").Value = "  -3.29%  "
$ws.Range("D18").Value = "'86.59"
$ws.Range("E18").Value = "  +5.06%  "
$ws.Range("D19").Value = "'0.05950"
$ws.Range("E19").Value = "  -11.46%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").Value = "'16.72"
$ws.Range("E21").Value = "  -3.88%  "
$ws.Range("D22").Value = "'6.062"
$ws.Range("E22").Value = "  -5.24%  "
$ws.Range("D23").Value = "'0.5256"
$ws.Range("E23").Value = "  -4.41%  "
$ws.Range("D24").Value = "'27.267.73"
$ws.Range("E24").Value = "  -3.02%  "
$ws.Range("D25").Value = "'11.29"
$ws.Range("E25").Value = "  -4.97%  "
$ws.Range("D26").Value = "'2.389"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").Value = "'19.98"
$ws.Range("E27").Value = "  -3.54%  "
$ws.Range("D28").Value = "'2.335"
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("D29").Value = "'148.04"
$ws.Range("E29").Value = "  -2.67%  "
$ws.Range("D30").Value = "'1.925.81"
$ws.Range("E30").Value = "  -3.79%  "
$ws.Range("D31").Value = "'125.60"
$ws.Range("E31").Value = "  -5.97%  "
$ws.Range("D32").Value = "'1.193"
$ws.Range("E32").Value = "  -5.15%  "
$ws.Range("D33").Value = "'0.09059"
$ws.Range("E33").Value = "  -5.96%  "
$ws.Range("D34").Value = "'5.573"
$ws.Range("E34").Value = "  -5.32%  "
$ws.Range("D35").Value = "'3.562"
$ws.Range("E35").Value = "  -12.16%  "
$ws.Range("D36").Value = "'12.67"
$ws.Range("E36").Value = "  +4.37%  "
$ws.Range("D37").Value = "'0.2134"
$ws.Range("E37").Value = "  -3.84%  "
$ws.Range("D38").Value = "'5.046"
$ws.Range("E38").Value = "  -4.00%  "
$ws.Range("D39").Value = "'0.06035"
$ws.Range("E39").Value = "  -4.73%  "
$ws.Range("D40").Value = "'0.02233"
$ws.Range("E40").Value = "  -6.03%  "
$ws.Range("D41").Value = "'0.6331"
$ws.Range("E41").Value = "  -5.55%  "
$ws.Range("D42").Value = "'1.184"
$ws.Range("E42").Value = "  -4.17%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "'1.407"
$ws.Range("E44").Value = "  -6.37%  "
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("D46").Value = "'13.39"
$ws.Range("E46").Value = "  -5.19%  "
$ws.Range("D47").Value = "'3.718"
$ws.Range("E47").Value = "  -3.31%  "
$ws.Range("D48").Value = "'0.5787"
$ws.Range("E48").Value = "  -5.99%  "
$ws.Range("D49").Value = "'123.82"
$ws.Range("E49").Value = "  -4.99%  "
$ws.Range("D50").Value = "'1.939"
$ws.Range("E50").Value = "  -5.49%  "
$ws.Range("D51").Value = "'0.06803"
$ws.Range("E51").Value = "  -4.56%  "
